$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column cells that already look numeric stay as text,
# matching the workbook convention where column D stores prices as strings.
$priceCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values (scraped symbol list refresh).
$ws.Range("D2").Value = '237.38'
$ws.Range("D3").Value = '22.01'
$ws.Range("D4").Value = '5.357'
$ws.Range("D5").Value = '0.05589'
$ws.Range("D6").Value = '6.479'
$ws.Range("D7").Value = '3.339'
$ws.Range("D8").Value = '0.7987'
$ws.Range("D9").Value = '1.045'
$ws.Range("D10").Value = '0.1386'
$ws.Range("D11").Value = '0.07302'
$ws.Range("D12").Value = '0.03137'
$ws.Range("D13").Value = '0.02953'
$ws.Range("D14").Value = '0.09233'
$ws.Range("D15").Value = '0.001665'
$ws.Range("D16").Value = '3.250'
$ws.Range("D17").Value = '0.04782'
$ws.Range("D18").Value = '0.0005714'
$ws.Range("E18").Value = '17OneONEWorstin24h'
$ws.Range("D19").Value = '0.006255'
$ws.Range("D20").Value = '0.005061'
$ws.Range("D21").Value = '0.001052'
$ws.Range("D22").Value = '0.0001500'
$ws.Range("D23").Value = '0.0003701'
$ws.Range("D24").Value = '3.973'
$ws.Range("D25").Value = '2.201'
$ws.Range("D40").Value = '0.04110'
$ws.Range("D41").Value = '0.007033'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").Value = '0.003501'
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("B43").Value = 'BKEXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D43").Value = '0.1038'
$ws.Range("E43").Value = '42BKEXTokenBKK'
$ws.Range("D44").Value = '0.008795'
$ws.Range("D45").Value = '0.00005440'
$ws.Range("D46").Value = '0.00000000750'
$ws.Range("D47").Value = '0.6756'
$ws.Range("D48").Value = '0.03630'
$ws.Range("E48").Value = '47BOLOBOLO'
$ws.Range("D49").Value = '0.00002101'
$ws.Range("D50").Value = '0.01010'
